$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (index 1): A1:I12 -> A1:I13
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Update "想去人数" (interested-count) figures that ticked up.
$ws1.Cells.Item(2, 6).Value = 62
$ws1.Cells.Item(3, 6).Value = 66
$ws1.Cells.Item(4, 6).Value = 161
$ws1.Cells.Item(5, 6).Value = 356
$ws1.Cells.Item(6, 6).Value = 5301
$ws1.Cells.Item(8, 6).Value = 5358
$ws1.Cells.Item(9, 6).Value = 621
$ws1.Cells.Item(11, 6).Value = 1369

# Shift the last existing row (old row 12, 九江 event) down to row 13,
# copying the whole range so value types/text (incl. the date-looking
# text in column B) and styles survive untouched.
$ws1.Range("A12:I12").Copy($ws1.Range("A13:I13"))
$ws1.Cells.Item(13, 1).Value = 12
$ws1.Cells.Item(13, 6).Value = 106

# Row 12 becomes the newly added 赣州 event; the index cell (A12) keeps
# its existing value (11) and style, only B12:I12 change. B11 already
# holds the literal text "2025-01-01" we need for B12 - copy the cell
# (instead of a Value round-trip) so the text isn't re-parsed as a date.
$ws1.Range("B11").Copy($ws1.Range("B12"))
$ws1.Cells.Item(12, 3).Value = "赣州·CA01动漫游戏嘉年华"
$ws1.Cells.Item(12, 4).Value = "廉泉路赣友味餐厅旁 铸谊篮球·羽毛球馆"
$ws1.Cells.Item(12, 5).Value = "2025.01.01 10:00-01.01 17:00"
$ws1.Cells.Item(12, 6).Value = 1
$ws1.Cells.Item(12, 7).Value = 19.9
$ws1.Cells.Item(12, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93787"
$ws1.Cells.Item(12, 9).Value = "//i1.hdslb.com/bfs/openplatform/202410/ZjFMZOiY1729500245727.jpeg"

# ---------------------------------------------------------------------------
# Sheet "全部类型" (index 4): A1:I13 -> A1:I14
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(2, 6).Value = 62
$ws4.Cells.Item(3, 6).Value = 66
$ws4.Cells.Item(4, 6).Value = 161
$ws4.Cells.Item(6, 6).Value = 356
$ws4.Cells.Item(7, 6).Value = 5301
$ws4.Cells.Item(9, 6).Value = 5358
$ws4.Cells.Item(10, 6).Value = 621
$ws4.Cells.Item(12, 6).Value = 1369

$ws4.Range("A13:I13").Copy($ws4.Range("A14:I14"))
$ws4.Cells.Item(14, 1).Value = 13
$ws4.Cells.Item(14, 6).Value = 106

$ws4.Range("B12").Copy($ws4.Range("B13"))
$ws4.Cells.Item(13, 3).Value = "赣州·CA01动漫游戏嘉年华"
$ws4.Cells.Item(13, 4).Value = "廉泉路赣友味餐厅旁 铸谊篮球·羽毛球馆"
$ws4.Cells.Item(13, 5).Value = "2025.01.01 10:00-01.01 17:00"
$ws4.Cells.Item(13, 6).Value = 1
$ws4.Cells.Item(13, 7).Value = 19.9
$ws4.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93787"
$ws4.Cells.Item(13, 9).Value = "//i1.hdslb.com/bfs/openplatform/202410/ZjFMZOiY1729500245727.jpeg"
